$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New / moved query text (shared strings content)
$qB2 = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE p.program_acronym IN ["TAILORx"]  
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`  order By ss.study_subject_id ASC LIMIT 100 
'@
$qC2 = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)  
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
WHERE p.program_acronym IN ["TAILORx"]
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@
$qB4 = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['Beagle'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS `Case ID` , coalesce(s.clinical_study_designation,'') AS `Study Code` , coalesce(s.clinical_study_type,'') AS  `Study Type`, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  `Neutered Status`
 order By ss.study_subject_id ASC LIMIT 100
'@
$qC4 = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Beagle']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

# Row 2: replace the two query cells with the new queries
$ws.Range("B2").Value = $qB2
$ws.Range("C2").Value = $qC2
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

# Row 4 (new row): old queries moved here
$ws.Range("B4").Value = $qB4
$ws.Range("C4").Value = $qC4
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# Row heights
$ws.Rows.Item(2).RowHeight = 316.8
$ws.Rows.Item(4).RowHeight = 244.8

# Column C width (was sized for the long query text, now much narrower)
$ws.Columns.Item(3).ColumnWidth = 47.166666666666664

# View: zoom + selection, matching the saved workbook view
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("B2").Select()
